# Applies the "Automatic update of files" change:
#  - Column C (Förändrad) for rows 2..11 bumped from 46062 to 46063 (i.e. +1 day)
#  - Rows 7..11 (columns A, B, G) rotated: the last row (11) moves to the top of
#    that block (row 7), and the others shift down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: bump column C (Förändrad) for data rows 2..11 by one day ---
for ($r = 2; $r -le 11; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = $cell.Value2 + 1
}

# --- Step 2: rotate rows 7..11 for columns A (1), B (2), G (7) ---
$colsToRotate = @(1, 2, 7)

foreach ($col in $colsToRotate) {
    # Capture current values for rows 7..11 in this column
    $values = @{}
    for ($r = 7; $r -le 11; $r++) {
        $values[$r] = $ws.Cells.Item($r, $col).Value2
    }

    # Row 7 gets old row 11's value; rows 8..11 get old row(r-1)'s value
    $ws.Cells.Item(7, $col).Value = $values[11]
    for ($r = 8; $r -le 11; $r++) {
        $ws.Cells.Item($r, $col).Value = $values[$r - 1]
    }
}
